# Inventory - assembled boards, updated boms
#
# Marks the common passive/jumper/LED "Status" column cells as USED
# (green fill) across the power-linear and power-switching-TPS62111 BOM
# sheets, and appends a revision-history row documenting the assembled
# 1x board, dated 2019-10-02 (serial 43740).

$wb = $excel.ActiveWorkbook

$usedGreen = 5296274   # BGR for RGB(0x92,0xD0,0x50) == FF92D050
$dateFmt   = "yyyy\-mm\-dd;@"
$centerAlign = -4108   # xlCenter

function Mark-Used($ws, [int[]]$rows) {
    foreach ($row in $rows) {
        $cell = $ws.Range("B$row")
        $cell.Value = "USED"
        $cell.Interior.Color = $usedGreen
    }
}

function Add-HistoryRow($ws, [int]$row, [int]$dateSerial, [string]$note) {
    $dateCell = $ws.Range("B$row")
    $dateCell.Value = $dateSerial
    $dateCell.NumberFormat = $dateFmt
    $dateCell.HorizontalAlignment = $centerAlign

    $ws.Range("C$row").Value = $note
}

# --- Sheet 1: power-linear-LP38692-2.7V-adj ---------------------------
$ws1 = $wb.Worksheets.Item(1)
Mark-Used $ws1 @(8,9,10,11,12,13,14,16,17,18,19)
Add-HistoryRow $ws1 30 43740 "Assembled 1x board"
[void]$ws1.Range("B10").Select()

# --- Sheet 2: power-linear-LP38692-3.3V-fix ----------------------------
$ws2 = $wb.Worksheets.Item(2)
Mark-Used $ws2 @(8,9,10,11,12,13,14,16,17)
Add-HistoryRow $ws2 28 43740 "Assembled 1x board"
[void]$ws2.Range("B10").Select()

# --- Sheet 3: power-linear-LP38692-5.0V-fix ----------------------------
$ws3 = $wb.Worksheets.Item(3)
[void]$ws3.Range("B20").Select()

# --- Sheet 5: power-switching-TPS62112-5.0 -----------------------------
# Selected/activated before sheet 4 so that sheet 4 ends up as the
# workbook's active tab (matches the tabSelected swap in the diff).
$ws5 = $wb.Worksheets.Item(5)
[void]$ws5.Range("B25").Select()

# --- Sheet 4: power-switching-TPS62111-3.3V ----------------------------
$ws4 = $wb.Worksheets.Item(4)
Mark-Used $ws4 @(11,12,13,14,15,16,17,19,20,21,22)
[void]$ws4.Range("B24").Select()
